$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("H1").Value = "Height"
$ws.Range("I1").Value = "Weight"

# New data rows
$ws.Range("H2").Value = -170
$ws.Range("I2").Value = 45

$ws.Range("H3").Value = 160
$ws.Range("I3").Value = 0

$ws.Range("H4").Value = 168468
$ws.Range("I4").Value = 54684

$ws.Range("H5").Value = 180
$ws.Range("I5").Value = 84

$ws.Range("H6").Value = 150.9
$ws.Range("I6").Value = 56.4

# Update selection
$ws.Range("K4").Select()

# Update window view (position/size of the workbook window)
$win = $excel.ActiveWindow
$win.Left = -108
$win.Top = -108
$win.Width = 23256
$win.Height = 12456
